# Update the metrics table on the active sheet:
#  - Row 2 (MAE) values updated
#  - Row 3 (MSE) values updated
#  - Row 4 label changed from "mean Y-Test" to "R2", values replaced with
#    new R2 scores
#  - Old rows 5 ("mean Y-predicted") and 6 ("R2") removed entirely, shrinking
#    the used range from A1:F6 to A1:F4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: MAE
$ws.Range("B2").Value = 0.994
$ws.Range("C2").Value = 1.05
$ws.Range("D2").Value = 1.157
$ws.Range("E2").Value = 0.5600000000000001
$ws.Range("F2").Value = 1.352

# Row 3: MSE
$ws.Range("B3").Value = 1.515
$ws.Range("C3").Value = 2.115
$ws.Range("D3").Value = 2.511
$ws.Range("E3").Value = 0.541
$ws.Range("F3").Value = 3.049

# Row 4: relabel to R2 and set new values
$ws.Range("A4").Value = "R2"
$ws.Range("B4").Value = 0.761
$ws.Range("C4").Value = 0.801
$ws.Range("D4").Value = 0.877
$ws.Range("E4").Value = 0.729
$ws.Range("F4").Value = 0.894

# Remove the old "mean Y-predicted" (row 5) and "R2" (row 6) rows so the
# sheet's used range shrinks to A1:F4, shifting any rows below them up.
$ws.Range("A5:F6").Delete()
